# "Zeitblatter auf stand gebracht" - bring Ursus Schneider's October
# timesheet up to date: log two more work entries (rows 26/27) and move
# the selection to reflect where editing left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26: 1h, DMX/ESP further work, 11:00 - 12:00
$ws.Range("B26").Value = 1
$ws.Range("C26").Value = "DMX und ESP weitere arbeit"
$ws.Range("D26").Value = "11:00 - 12:00"

# Row 27: 2h, DMX/ESP finished, 14:00 - 16:00
$ws.Range("B27").Value = 2
$ws.Range("C27").Value = "DMX und ESP ferigstellen"
$ws.Range("D27").Value = "14:00 - 16:00"

# Leave the selection where the author ended up editing.
$ws.Range("D28").Select()
